$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of results for consumer/lame, mirroring the style of the
# existing "indented" rows (A3:A5) which use the bold-ish style (s="1").
$ws.Range("A7").Value = "consumer/lame/runme_large.sh"
$ws.Range("B7").Value = 0.17
$ws.Range("C7").Value = 0.16
$ws.Range("D7").Value = 0

# Match the formatting used by the other "sub-item" rows (A3:A5).
$ws.Range("A7").Style = $ws.Range("A3").Style

# Move the selection to the newly added row, as in the authored workbook.
$ws.Range("A7").Select()
